$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change), including three
# rows whose ranking order swapped (rows 19/20, 33/34, 39/40 -> new coin
# name/link at that rank together with its price/volume).
$updates = @(
    @{Row=2; Col='D'; Value='69.299.88'},
    @{Row=2; Col='E'; Value='  -2.09%  '},
    @{Row=3; Col='D'; Value='3.483.71'},
    @{Row=3; Col='E'; Value='  -3.46%  '},
    @{Row=4; Col='E'; Value='  +0.19%  '},
    @{Row=5; Col='D'; Value='576.54'},
    @{Row=5; Col='E'; Value='  -4.25%  '},
    @{Row=6; Col='D'; Value='189.29'},
    @{Row=6; Col='E'; Value='  -6.67%  '},
    @{Row=7; Col='D'; Value='0.616'},
    @{Row=7; Col='E'; Value='  -2.05%  '},
    @{Row=8; Col='D'; Value='3.475.37'},
    @{Row=8; Col='E'; Value='  -3.34%  '},
    @{Row=9; Col='E'; Value='  +0.07%  '},
    @{Row=10; Col='D'; Value='0.202'},
    @{Row=10; Col='E'; Value='  -6.58%  '},
    @{Row=11; Col='D'; Value='0.617'},
    @{Row=11; Col='E'; Value='  -4.84%  '},
    @{Row=12; Col='D'; Value='50.28'},
    @{Row=12; Col='E'; Value='  -6.13%  '},
    @{Row=13; Col='D'; Value='0.0000281'},
    @{Row=13; Col='E'; Value='  -6.89%  '},
    @{Row=14; Col='D'; Value='9.05'},
    @{Row=14; Col='E'; Value='  -5.87%  '},
    @{Row=15; Col='D'; Value='4.045.01'},
    @{Row=15; Col='E'; Value='  -3.25%  '},
    @{Row=16; Col='D'; Value='638.93'},
    @{Row=16; Col='E'; Value='  -6.72%  '},
    @{Row=17; Col='D'; Value='69.155.75'},
    @{Row=17; Col='E'; Value='  -2.41%  '},
    @{Row=18; Col='D'; Value='3.496.88'},
    @{Row=18; Col='E'; Value='  -3.05%  '},
    @{Row=19; Col='B'; Value='TRON'},
    @{Row=19; Col='C'; Value='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'},
    @{Row=19; Col='D'; Value='0.120'},
    @{Row=19; Col='E'; Value='  -2.38%  '},
    @{Row=20; Col='B'; Value='Uniswap'},
    @{Row=20; Col='C'; Value='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'},
    @{Row=20; Col='D'; Value='12.21'},
    @{Row=20; Col='E'; Value='  -4.30%  '},
    @{Row=21; Col='D'; Value='18.24'},
    @{Row=21; Col='E'; Value='  -4.83%  '},
    @{Row=22; Col='D'; Value='0.942'},
    @{Row=22; Col='E'; Value='  -5.68%  '},
    @{Row=23; Col='D'; Value='18.07'},
    @{Row=23; Col='E'; Value='  -2.97%  '},
    @{Row=24; Col='D'; Value='5.22'},
    @{Row=24; Col='E'; Value='  -3.45%  '},
    @{Row=25; Col='D'; Value='98.52'},
    @{Row=25; Col='E'; Value='  -10.99%  '},
    @{Row=26; Col='D'; Value='4.25'},
    @{Row=26; Col='E'; Value='  -7.98%  '},
    @{Row=27; Col='D'; Value='2.85'},
    @{Row=27; Col='E'; Value='  -5.89%  '},
    @{Row=28; Col='D'; Value='9.84'},
    @{Row=28; Col='E'; Value='  -7.20%  '},
    @{Row=29; Col='D'; Value='9.24'},
    @{Row=29; Col='E'; Value='  -9.79%  '},
    @{Row=30; Col='D'; Value='32.17'},
    @{Row=30; Col='E'; Value='  -6.84%  '},
    @{Row=31; Col='D'; Value='4.09'},
    @{Row=31; Col='E'; Value='  -10.07%  '},
    @{Row=32; Col='D'; Value='6.60'},
    @{Row=32; Col='E'; Value='  -8.63%  '},
    @{Row=33; Col='B'; Value='Cosmos'},
    @{Row=33; Col='C'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Row=33; Col='D'; Value='11.48'},
    @{Row=33; Col='E'; Value='  -5.99%  '},
    @{Row=34; Col='B'; Value='Bittensor'},
    @{Row=34; Col='C'; Value='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'},
    @{Row=34; Col='D'; Value='573.40'},
    @{Row=34; Col='E'; Value='  +12.64%  '},
    @{Row=35; Col='D'; Value='0.108'},
    @{Row=35; Col='E'; Value='  -5.84%  '},
    @{Row=36; Col='D'; Value='60.36'},
    @{Row=36; Col='E'; Value='  -5.18%  '},
    @{Row=37; Col='D'; Value='3.844.51'},
    @{Row=37; Col='E'; Value='  -1.41%  '},
    @{Row=38; Col='E'; Value='  -0.07%  '},
    @{Row=39; Col='B'; Value='CoreDAO'},
    @{Row=39; Col='C'; Value='https://coinranking.com/coin/HFvoXUQh4+coredao-core'},
    @{Row=39; Col='D'; Value='3.80'},
    @{Row=39; Col='E'; Value='  +35.58%  '},
    @{Row=40; Col='B'; Value='PEPE'},
    @{Row=40; Col='C'; Value='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'},
    @{Row=40; Col='D'; Value='0.0₃0778'},
    @{Row=40; Col='E'; Value='  -8.56%  '},
    @{Row=41; Col='D'; Value='3.49'},
    @{Row=41; Col='E'; Value='  -2.65%  '},
    @{Row=42; Col='D'; Value='2.84'},
    @{Row=42; Col='E'; Value='  -6.23%  '},
    @{Row=43; Col='D'; Value='0.368'},
    @{Row=43; Col='E'; Value='  -4.76%  '},
    @{Row=44; Col='E'; Value='  -5.72%  '},
    @{Row=45; Col='D'; Value='33.26'},
    @{Row=45; Col='E'; Value='  -10.33%  '},
    @{Row=46; Col='D'; Value='0.0437'},
    @{Row=46; Col='E'; Value='  -6.70%  '},
    @{Row=47; Col='D'; Value='3.33'},
    @{Row=47; Col='E'; Value='  -2.65%  '},
    @{Row=48; Col='D'; Value='2.81'},
    @{Row=48; Col='E'; Value='  -8.43%  '},
    @{Row=49; Col='E'; Value='  -4.30%  '},
    @{Row=50; Col='D'; Value='0.999'},
    @{Row=50; Col='E'; Value='  -0.25%  '},
    @{Row=51; Col='D'; Value='8.12'},
    @{Row=51; Col='E'; Value='  -6.06%  '}
)

foreach ($item in $updates) {
    $addr = "$($item.Col)$($item.Row)"
    $cell = $ws.Range($addr)
    if ($item.Col -eq "D") {
        # Force text storage so numeric-looking prices ("576.54", "0.616", ...)
        # are written exactly as the original inline-string text instead of
        # being auto-converted to floating point numbers by Excel. Then drop
        # the temporary text number-format again so the cell's style index
        # matches the untouched (default/no-format) cells around it.
        $cell.NumberFormat = "@"
        $cell.Value = $item.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $item.Value
    }
}
